# Add a dungeon card shop: append "|cardshop;2" to the QuestDungeonRate
# strings used by the two dungeon rows whose M column currently reads
# "fight;7|..." (row 4) and "fight;10|..." (row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M4").Value = $ws.Range("M4").Value2 + "|cardshop;2"
$ws.Range("M6").Value = $ws.Range("M6").Value2 + "|cardshop;2"
